$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row swap: Hedera (34) <-> InternetComputer(DFINITY) (35),
# and Aave (42) <-> InjectiveProtocol (43), plus per-row price/volume updates.

$ws.Range("D2").Value = "37.934.26"
$ws.Range("E2").Value = "  -0.34%  "
$ws.Range("D3").Value = "2.039.20"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").Value = "'60.40"
$ws.Range("E7").Value = "  +4.37%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").Value = "'0.388"
$ws.Range("E11").Value = "  +0.89%  "
$ws.Range("D12").Value = "'14.68"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").Value = "2.339.30"
$ws.Range("E13").Value = "  -1.07%  "
$ws.Range("D14").Value = "'21.12"
$ws.Range("E14").Value = "  +2.38%  "
$ws.Range("D15").Value = "'0.764"
$ws.Range("E15").Value = "  +1.48%  "
$ws.Range("D16").Value = "'5.21"
$ws.Range("E16").Value = "  -1.44%  "
$ws.Range("D17").Value = "2.023.92"
$ws.Range("E17").Value = "  -2.19%  "
$ws.Range("D18").Value = "37.846.03"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").Value = "'6.10"
$ws.Range("E19").Value = "  -1.50%  "
$ws.Range("D20").Value = "'69.89"
$ws.Range("E20").Value = "  +0.23%  "
$ws.Range("D21").Value = "0.0₃0826"
$ws.Range("E21").Value = "  -0.45%  "
$ws.Range("D22").Value = "'225.45"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("E23").Value = "  -0.06%  "
$ws.Range("E24").Value = "  -2.41%  "
$ws.Range("E25").Value = "  -1.95%  "
$ws.Range("D26").Value = "'9.27"
$ws.Range("E26").Value = "  -0.42%  "
$ws.Range("D27").Value = "'165.24"
$ws.Range("E27").Value = "  -0.23%  "
$ws.Range("E28").Value = "  -4.05%  "
$ws.Range("D29").Value = "'18.95"
$ws.Range("E29").Value = "  -0.40%  "
$ws.Range("E30").Value = "  -5.97%  "
$ws.Range("E31").Value = "  +1.55%  "
$ws.Range("D32").Value = "'4.46"
$ws.Range("E32").Value = "  -2.04%  "
$ws.Range("D33").Value = "'2.06"
$ws.Range("E33").Value = "  +3.56%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "'4.51"
$ws.Range("E34").Value = "  -2.18%  "
$ws.Range("B35").Value = "Hedera"
$ws.Range("C35").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D35").Value = "'0.0602"
$ws.Range("E35").Value = "  -1.91%  "
$ws.Range("E36").Value = "  +6.68%  "
$ws.Range("E37").Value = "  -4.84%  "
$ws.Range("D38").Value = "'3.25"
$ws.Range("E38").Value = "  -2.02%  "
$ws.Range("D39").Value = "'0.999"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").Value = "1.542.28"
$ws.Range("E40").Value = "  +4.04%  "
$ws.Range("D41").Value = "'0.0218"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").Value = "'16.95"
$ws.Range("E42").Value = "  +0.96%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'97.12"
$ws.Range("E43").Value = "  -1.33%  "
$ws.Range("E44").Value = "  -1.27%  "
$ws.Range("D45").Value = "'0.0924"
$ws.Range("E45").Value = "  -2.15%  "
$ws.Range("E46").Value = "  -1.16%  "
$ws.Range("D47").Value = "'3.94"
$ws.Range("E47").Value = "  -3.82%  "
$ws.Range("E48").Value = "  -1.28%  "
$ws.Range("E49").Value = "  -0.37%  "
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("D51").Value = "2.228.22"
$ws.Range("E51").Value = "  -0.96%  "
